$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.390.20'
$ws.Range('E2').Value = '  +1.94%  '
$ws.Range('D3').Value = '1.842.27'
$ws.Range('E3').Value = '  +1.51%  '
$ws.Range('D4').Value = '1.015'
$ws.Range('E4').Value = '  +1.40%  '
$ws.Range('D5').Value = '315.54'
$ws.Range('E5').Value = '  +2.21%  '
$ws.Range('D6').Value = '1.013'
$ws.Range('E6').Value = '  +1.20%  '
$ws.Range('D7').Value = '0.4747'
$ws.Range('E7').Value = '  +1.68%  '
$ws.Range('D8').Value = '0.3702'
$ws.Range('D9').Value = '0.07476'
$ws.Range('E9').Value = '  +1.36%  '
$ws.Range('D10').Value = '0.8862'
$ws.Range('E10').Value = '  +1.88%  '
$ws.Range('D11').Value = '20.53'
$ws.Range('E11').Value = '  +0.71%  '
$ws.Range('D12').Value = '1.845.09'
$ws.Range('E12').Value = '  +2.82%  '
$ws.Range('D13').Value = '0.07376'
$ws.Range('E13').Value = '  +4.42%  '
$ws.Range('D14').Value = '5.491'
$ws.Range('E14').Value = '  +2.56%  '
$ws.Range('D15').Value = '93.34'
$ws.Range('E15').Value = '  +1.91%  '
$ws.Range('D16').Value = '6.585'
$ws.Range('E16').Value = '  +1.40%  '
$ws.Range('E17').Value = '  +1.29%  '
$ws.Range('D18').Value = '0.000008865'
$ws.Range('E18').Value = '  +1.95%  '
$ws.Range('E19').Value = '  +1.23%  '
$ws.Range('D20').Value = '14.86'
$ws.Range('E20').Value = '  +1.02%  '
$ws.Range('D21').Value = '27.415.21'
$ws.Range('E21').Value = '  +1.97%  '
$ws.Range('D22').Value = '5.351'
$ws.Range('E22').Value = '  +0.43%  '
$ws.Range('E23').Value = '  +1.51%  '
$ws.Range('D24').Value = '2.079.71'
$ws.Range('E24').Value = '  +0.94%  '
$ws.Range('E25').Value = '  +0.39%  '
$ws.Range('D26').Value = '152.20'
$ws.Range('E26').Value = '  +1.59%  '
$ws.Range('D27').Value = '18.67'
$ws.Range('E27').Value = '  +1.97%  '
$ws.Range('D28').Value = '2.183'
$ws.Range('E28').Value = '  +0.64%  '
$ws.Range('D29').Value = '5.278'
$ws.Range('E29').Value = '  -0.76%  '
$ws.Range('D30').Value = '118.10'
$ws.Range('E30').Value = '  +2.11%  '
$ws.Range('D31').Value = '0.08971'
$ws.Range('E31').Value = '  +0.45%  '
$ws.Range('D32').Value = '0.7623'
$ws.Range('E32').Value = '  -0.53%  '
$ws.Range('E33').Value = '  +1.61%  '
$ws.Range('D34').Value = '4.568'
$ws.Range('E34').Value = '  +1.44%  '
$ws.Range('D35').Value = '2.947'
$ws.Range('E35').Value = '  +1.59%  '
$ws.Range('E36').Value = '  +1.31%  '
$ws.Range('D37').Value = '1.106'
$ws.Range('E37').Value = '  +1.84%  '
$ws.Range('D38').Value = '0.05371'
$ws.Range('E38').Value = '  +1.74%  '
$ws.Range('D39').Value = '0.01963'
$ws.Range('E39').Value = '  +0.16%  '
$ws.Range('D40').Value = '3.002'
$ws.Range('E40').Value = '  +2.14%  '
$ws.Range('D41').Value = '7.329'
$ws.Range('E41').Value = '  +1.09%  '
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').Value = '0.5359'
$ws.Range('E42').Value = '  +0.91%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').Value = '2.377'
$ws.Range('E43').Value = '  +1.06%  '
$ws.Range('D44').Value = '0.1670'
$ws.Range('E44').Value = '  +0.52%  '
$ws.Range('D45').Value = '8.551'
$ws.Range('E45').Value = '  +1.65%  '
$ws.Range('D46').Value = '0.4982'
$ws.Range('E46').Value = '  +1.17%  '
$ws.Range('D47').Value = '10.52'
$ws.Range('E47').Value = '  +0.49%  '
$ws.Range('E48').Value = '  +1.34%  '
$ws.Range('D49').Value = '105.19'
$ws.Range('E49').Value = '  +1.54%  '
$ws.Range('D50').Value = '1.683'
$ws.Range('E50').Value = '  +0.86%  '
$ws.Range('D51').Value = '0.06335'
$ws.Range('E51').Value = '  +0.81%  '
